# fix: changement de calcul de l'évaluation du gain + corrections de la doc
# Rename the "Gain" column of Tableau1 to "Evolution" and change its formula
# from (avec/sans) to ((avec-sans)/sans) for both the per-row calculated
# column and the totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculatrice")
$lo = $ws.ListObjects.Item("Tableau1")

# Rename the table column by editing its header cell (this keeps the table's
# structured references, dxf/style links, etc. in sync).
$headerCell = $ws.Range("R3")
$headerCell.Value = "Evolution"

# Update the calculated formula for the whole data column.
$col = $lo.ListColumns.Item("Evolution")
$col.DataBodyRange.Formula = "=(Tableau1[[#This Row],[kgCO2eq avec kube-downscaler]]-Tableau1[[#This Row],[kgCO2eq sans kube-downscaler]])/Tableau1[[#This Row],[kgCO2eq sans kube-downscaler]]"

# Update the totals row formula for that column.
$totalsCell = $ws.Range("R9")
$totalsCell.Formula = "=(Tableau1[[#Totals],[kgCO2eq avec kube-downscaler]]-Tableau1[[#Totals],[kgCO2eq sans kube-downscaler]])/Tableau1[[#Totals],[kgCO2eq sans kube-downscaler]]"

# Reflect the new zoom/selection state on the "Calculatrice" sheet view.
$ws.Activate()
$excel.ActiveWindow.Zoom = 93
$excel.ActiveWindow.View = $excel.ActiveWindow.View
$ws.Range("Q12").Select()

$wb.Save()
